$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44342

# Row 4
$ws.Range("D4").Value = 44372

# Row 5
$ws.Range("D5").Value = 44354
$ws.Range("O5").Value = "Región del Maule"

# Row 6
$ws.Range("D6").Value = 44371
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 6500
$ws.Range("M6").Value = 6500
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 181

# Row 7
$ws.Range("D7").Value = 44348

# Row 8
$ws.Range("D8").Value = 44386
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 6500
$ws.Range("M8").Value = 6500
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 181

# Row 9
$ws.Range("D9").Value = 44376
$ws.Range("J9").Value = 150

# Row 10
$ws.Range("D10").Value = 44355
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 7000
$ws.Range("P10").Value = 194

# Row 11
$ws.Range("D11").Value = 44340

# Row 12
$ws.Range("D12").Value = 44369
$ws.Range("N12").Value = "`$/caja 20 docenas"
$ws.Range("P12").Value = 7000
$ws.Range("Q12").Value = 1

# Row 13
$ws.Range("D13").Value = 44690
$ws.Range("J13").Value = 500
$ws.Range("O13").Value = "Región del Maule"

# Row 14
$ws.Range("D14").Value = 44357
$ws.Range("N14").Value = "`$/caja 20 docenas"
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 6500
$ws.Range("Q14").Value = 1

# Row 15
$ws.Range("D15").Value = 44364
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = "`$/caja 36 atados"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 194
$ws.Range("Q15").Value = 36

# Row 16
$ws.Range("D16").Value = 44358
$ws.Range("J16").Value = 150
$ws.Range("N16").Value = "`$/caja 36 atados"
$ws.Range("P16").Value = 194
$ws.Range("Q16").Value = 36
